$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based index of the first paragraph whose (trimmed) text
# equals $needle.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($needle, $startAt) {
    for ($i = $startAt; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $needle) {
            return $i
        }
    }
    return -1
}

$titleText = "Play Astro Pug for Free: Chinese Tradition Meets Pugs"
$oldMetaText = "Get ready to play Astro Pug for free! With a unique combination of Chinese architecture and pugs, this game features exciting free spins and multipliers."
$newImagePrompt = "Create a feature image for Astro Pug: Design a cartoon-style image featuring a happy Maya warrior with glasses playing Astro Pug. The warrior should be sitting in front of the slot machine with a big smile on his face and hand gestures showing excitement. He should be wearing traditional Maya clothing, including a headpiece and necklace, and have glasses on. The background should feature Chinese architecture structures with red roofs and blue skies. The Astro Pug slot machine should be in the foreground, with the reels spinning and the Pug symbol prominently displayed. The image should be colorful, vibrant, and eye-catching to attract viewers' attention."

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titleIdx = Find-ParagraphIndex $titleText 1

$titleEnd = $d.Paragraphs($titleIdx).Range
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

$metaIdx = $titleIdx + 1
$metaPara = $d.Paragraphs($metaIdx)
$metaPara.Style = "Normal"

$metaRange = $d.Paragraphs($metaIdx).Range
$metaRange.Collapse(1)
$metaStart = $metaRange.Start
$metaRange.InsertBefore("Meta description" + ": " + $oldMetaText)

# Bold just the "Meta description" label (16 characters), leaving the rest
# (including the leading colon) in regular weight.
$boldRange = $d.Range($metaStart, $metaStart + 16)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold title
#    paragraph and replace the italic meta-description paragraph's text
#    with the new feature-image prompt (keeping its italic formatting).
# ---------------------------------------------------------------------------
$searchFrom = $metaIdx + 1
$dupTitleIdx = Find-ParagraphIndex $titleText $searchFrom
if ($dupTitleIdx -ne -1) {
    $d.Paragraphs($dupTitleIdx).Range.Delete()
}

$oldMetaIdx = Find-ParagraphIndex $oldMetaText $searchFrom
$oldMetaPara = $d.Paragraphs($oldMetaIdx)
$oldMetaRange = $d.Range($oldMetaPara.Range.Start, $oldMetaPara.Range.End)
$oldMetaRange.Text = $newImagePrompt

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
Write-Output ("Paragraph " + $metaIdx + ": " + $d.Paragraphs($metaIdx).Range.Text)
Write-Output ("Paragraph " + $oldMetaIdx + ": " + $d.Paragraphs($oldMetaIdx).Range.Text)
